$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 79245
$ws.Range("A4").Value = 131256691
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 488667
$ws.Range("R4").Value = 6665262
$ws.Range("Z4").Value = "10:55"
$ws.Range("AB4").Value = "10:55"
$ws.Range("AC4").Value = "Ringhack på gran."
$ws.Range("A5").Value = 131260583
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 488834
$ws.Range("R5").Value = 6665228
$ws.Range("Z5").Value = "15:30"
$ws.Range("AB5").Value = "15:30"
$ws.Range("AC5").Value = "Ringhack på tall."
$ws.Range("B6").Value = 79245
$ws.Range("B7").Value = 91830
$ws.Range("B8").Value = 91830
$ws.Range("B10").Value = 79245
$ws.Range("B11").Value = 79245
$ws.Range("A12").Value = 131260641
$ws.Range("Q12").Value = 488859
$ws.Range("R12").Value = 6665292
$ws.Range("Z12").Value = "15:34"
$ws.Range("AB12").Value = "15:34"
$ws.Range("AC12").Value = "Ringhack på gran."
$ws.Range("A13").Value = 131257290
$ws.Range("Q13").Value = 488842
$ws.Range("R13").Value = 6665224
$ws.Range("Z13").Value = "11:26"
$ws.Range("AB13").Value = "11:26"
$ws.Range("AC13").Value = "Ringhack på tall."
$ws.Range("A14").Value = 131256673
$ws.Range("Q14").Value = 488652
$ws.Range("R14").Value = 6665282
$ws.Range("Z14").Value = "10:54"
$ws.Range("AB14").Value = "10:54"
$ws.Range("B15").Value = 79245
$ws.Range("B16").Value = 91830
$ws.Range("B17").Value = 79245
$ws.Range("B18").Value = 79245
$ws.Range("B19").Value = 79245
$ws.Range("B20").Value = 79245
$ws.Range("B21").Value = 81230
$ws.Range("B22").Value = 79245
$ws.Range("B23").Value = 79245
$ws.Range("A24").Value = 131257045
$ws.Range("B24").Value = 79245
$ws.Range("Q24").Value = 488760
$ws.Range("R24").Value = 6665302
$ws.Range("Z24").Value = "11:10"
$ws.Range("AB24").Value = "11:10"
$ws.Range("AC24").Value = "Gran."
$ws.Range("A25").Value = 131257650
$ws.Range("B25").Value = 79245
$ws.Range("Q25").Value = 488911
$ws.Range("R25").Value = 6665227
$ws.Range("Z25").Value = "12:00"
$ws.Range("AB25").Value = "12:00"
$ws.Range("AC25").Value = "Gran"
$ws.Range("B26").Value = 79245
$ws.Range("A27").Value = 131255910
$ws.Range("B27").Value = 79245
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("M27").Value = ""
$ws.Range("Q27").Value = 488763
$ws.Range("R27").Value = 6665157
$ws.Range("Z27").Value = "10:03"
$ws.Range("AB27").Value = "10:03"
$ws.Range("AC27").Value = "Tall."
$ws.Range("A28").Value = 131258531
$ws.Range("B28").Value = 79245
$ws.Range("Q28").Value = 488725
$ws.Range("R28").Value = 6665212
$ws.Range("Z28").Value = "13:02"
$ws.Range("AB28").Value = "13:02"
$ws.Range("AC28").Value = "Gran"
$ws.Range("A29").Value = 131257239
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 488852
$ws.Range("R29").Value = 6665286
$ws.Range("Z29").Value = "11:23"
$ws.Range("AB29").Value = "11:23"
$ws.Range("AC29").Value = "Barkfläk, hagelsalva."
$ws.Range("A30").Value = 131258537
$ws.Range("B30").Value = 79245
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("M30").Value = ""
$ws.Range("Q30").Value = 488726
$ws.Range("R30").Value = 6665209
$ws.Range("Z30").Value = "13:02"
$ws.Range("AB30").Value = "13:02"
$ws.Range("AC30").Value = "Gran"
$ws.Range("B31").Value = 79245
$ws.Range("A32").Value = 131257659
$ws.Range("B32").Value = 57884
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = "Tretåig hackspett"
$ws.Range("G32").Value = "Picoides tridactylus"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("M32").Value = "äldre spår"
$ws.Range("Q32").Value = 488901
$ws.Range("R32").Value = 6665222
$ws.Range("Z32").Value = "12:02"
$ws.Range("AB32").Value = "12:02"
$ws.Range("AC32").Value = "Ringhack på tall."
$ws.Range("B33").Value = 79245
$ws.Range("B34").Value = 79245
$ws.Range("A35").Value = 131257385
$ws.Range("B35").Value = 91830
$ws.Range("E35").Value = 5432
$ws.Range("F35").Value = "Granticka"
$ws.Range("G35").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("Q35").Value = 488876
$ws.Range("R35").Value = 6665194
$ws.Range("Z35").Value = "11:31"
$ws.Range("AB35").Value = "11:31"
$ws.Range("AC35").Value = "Lågstubbe."
$ws.Range("A36").Value = 131260531
$ws.Range("B36").Value = 79245
$ws.Range("E36").Value = 6425
$ws.Range("F36").Value = "Garnlav"
$ws.Range("G36").Value = "Alectoria sarmentosa"
$ws.Range("H36").Value = "(Ach.) Ach."
$ws.Range("Q36").Value = 488786
$ws.Range("R36").Value = 6665188
$ws.Range("Z36").Value = "15:25"
$ws.Range("AB36").Value = "15:25"
$ws.Range("AC36").Value = "Gran"
$ws.Range("A37").Value = 131256459
$ws.Range("B37").Value = 57881
$ws.Range("E37").Value = 100049
$ws.Range("F37").Value = "Spillkråka"
$ws.Range("G37").Value = "Dryocopus martius"
$ws.Range("H37").Value = "(Linnaeus, 1758)"
$ws.Range("M37").Value = "färska spår"
$ws.Range("Q37").Value = 488669
$ws.Range("R37").Value = 6665268
$ws.Range("Z37").Value = "10:42"
$ws.Range("AB37").Value = "10:42"
$ws.Range("AC37").Value = "Färska och äldre hack."
$ws.Range("B38").Value = 79245
$ws.Range("B39").Value = 79245
$ws.Range("A40").Value = 131257035
$ws.Range("B40").Value = 79245
$ws.Range("E40").Value = 6425
$ws.Range("F40").Value = "Garnlav"
$ws.Range("G40").Value = "Alectoria sarmentosa"
$ws.Range("H40").Value = "(Ach.) Ach."
$ws.Range("M40").Value = ""
$ws.Range("Q40").Value = 488760
$ws.Range("R40").Value = 6665301
$ws.Range("Z40").Value = "11:09"
$ws.Range("AB40").Value = "11:09"
$ws.Range("AC40").Value = "Gran"
$ws.Range("A41").Value = 131257343
$ws.Range("B41").Value = 57884
$ws.Range("E41").Value = 100109
$ws.Range("F41").Value = "Tretåig hackspett"
$ws.Range("G41").Value = "Picoides tridactylus"
$ws.Range("H41").Value = "(Linnaeus, 1758)"
$ws.Range("M41").Value = "färska spår"
$ws.Range("Q41").Value = 488872
$ws.Range("R41").Value = 6665191
$ws.Range("Z41").Value = "11:29"
$ws.Range("AB41").Value = "11:29"
$ws.Range("AC41").Value = "Mycket barkfläk, hagelsalvor på många träd, skalade klena senvuxna granar."
$ws.Range("A42").Value = 131273991
$ws.Range("B42").Value = 79245
$ws.Range("Q42").Value = 488928
$ws.Range("R42").Value = 6665146
$ws.Range("A43").Value = 131273946
$ws.Range("B43").Value = 79245
$ws.Range("Q43").Value = 488774
$ws.Range("R43").Value = 6665353
$ws.Range("B45").Value = 79245
